$wb = $excel.ActiveWorkbook

# --- Astro sheet: B2 value change 200 -> 700 ---
$astro = $wb.Worksheets.Item("Astro")
$astro.Range("B2").Value = 700

# --- Insert new "FireSat" sheet between "Astro" and "AOCS" ---
$aocs = $wb.Worksheets.Item("AOCS")
$fireSat = $wb.Worksheets.Add($aocs)
$fireSat.Name = "FireSat"

$fireSat.Range("A1").Value = "name"
$fireSat.Range("B1").Value = "value"
$fireSat.Range("C1").Value = "units"

$fireSat.Range("A2").Value = "mass"
$fireSat.Range("B2").Value = 215

$fireSat.Range("A3").Value = "dipole"
$fireSat.Range("B3").Value = 0

$fireSat.Range("A4").Value = "cg"
$fireSat.Range("B4").Value = 0

$fireSat.Range("A5").Value = "c_pres aero"
$fireSat.Range("B5").Value = 0

$fireSat.Range("A6").Value = "c_pres solar"

$fireSat.Range("A7").Value = "solar incidence"

$fireSat.Range("A8").Value = "pt excursion"

$fireSat.Range("A9").Value = "Cd"

$fireSat.Range("A10").Value = "xdim"
$fireSat.Range("B10").Value = 1.294

$fireSat.Range("A11").Value = "ydim"
$fireSat.Range("B11").Value = 1.8299799999999999

$fireSat.Range("A12").Value = "zdim"
$fireSat.Range("B12").Value = 1.294

# Match the bold/centered/bordered header style used on the other sheets
$astro.Range("A1:C1").Copy() | Out-Null
$fireSat.Range("A1:C1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# size the columns to fit their content, like the other sheets in the workbook
$fireSat.Range("A1:C12").Columns.AutoFit() | Out-Null

# leave the selection where it was left in the authored file
$fireSat.Range("B6").Select() | Out-Null

# --- Re-select the intended cell on the Astro sheet, and make it the active tab ---
$astro.Activate()
$astro.Range("C14").Select() | Out-Null
